$d = $word.ActiveDocument

# --- Step 1: remove the heading paragraph entirely -------------------------
# Deleting the first paragraph's own Range (which, in the Word object model,
# includes its trailing paragraph mark) removes the "My first flower
# classifier with fast ai" text, drops the Heading1 style, and merges what
# used to be paragraph 1 with the (empty, non-heading) paragraph that
# followed it - the surviving paragraph keeps that second paragraph's mark
# (and therefore its plain, non-heading formatting).
$headingPara = $d.Paragraphs(1)
$headingPara.Range.Delete()

# --- Step 2: drop the "_GoBack" edit-location bookmark where it now belongs
# Word automatically re-homes the hidden "_GoBack" bookmark at the site of
# the most recent edit; re-adding a bookmark with that name moves it (Word
# bookmark names are unique, so this both removes the old occurrence, which
# sat between "Check this link" and the following space, and creates the
# new one in the now-empty first paragraph).
$firstPara = $d.Paragraphs(1)
$d.Bookmarks.Add("_GoBack", $firstPara.Range)

# --- Step 3: merge the "Check this link" / " " runs into a single run ------
# Re-finding and "replacing" the phrase (including the trailing space) with
# itself collapses the two adjacent same-formatted runs - and the bookmark
# that used to sit between them - into one run with text
# "Check this link ".
$d.Content.Find.Execute("Check this link ", $true, $false, $false, $false, $false, $true, 1, $false, "Check this link ", 2)
